$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "align setting" / "align execution" labels in column C for rows 24-25
# to the new "wafer setting" / "execution" labels.
$ws.Range("C24").Value = "웨이퍼 설정"
$ws.Range("C25").Value = "실행"

# Update the active selection on the sheet to E25 (previously D25).
$ws.Range("E25").Select()
